$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.2895605203232839
$ws.Range("J2").Value = 0.2895605203232838
$ws.Range("M2").Value = 51.15371566666666
$ws.Range("N2").Value = 153.461147
$ws.Range("O2").Value = 0.3311207986511828
$ws.Range("P2").Value = 0.3311207986511828
$ws.Range("Q2").Value = 69.76679652019544
$ws.Range("R2").Value = 627.9011686817589
$ws.Range("S2").Value = 0.09587951074729782
$ws.Range("T2").Value = 0.09587951074729779
# Row 3
$ws.Range("I3").Value = 0.2895605203232839
$ws.Range("J3").Value = 0.2895605203232838
$ws.Range("M3").Value = 53.36146666666667
$ws.Range("O3").Value = 0.3454116915964105
$ws.Range("P3").Value = 0.3454116915964106
$ws.Range("Q3").Value = 72.77787230964445
$ws.Range("R3").Value = 655.0008507868001
$ws.Range("S3").Value = 0.1000175891444023
$ws.Range("T3").Value = 0.1000175891444023
# Row 4
$ws.Range("I4").Value = 0.2895605203232839
$ws.Range("J4").Value = 0.2895605203232838
$ws.Range("M4").Value = 25.00653133333333
$ws.Range("N4").Value = 75.019594
$ws.Range("O4").Value = 0.1618686447050176
$ws.Range("P4").Value = 0.1618686447050176
$ws.Range("Q4").Value = 34.10554952795755
$ws.Range("R4").Value = 306.949945751618
$ws.Range("S4").Value = 0.04687076898480968
$ws.Range("T4").Value = 0.04687076898480967
# Row 5
$ws.Range("I5").Value = 0.2895605203232839
$ws.Range("J5").Value = 0.2895605203232838
$ws.Range("M5").Value = 6.481347
$ws.Range("N5").Value = 19.444041
$ws.Range("O5").Value = 0.0419541135381084
$ws.Range("P5").Value = 0.0419541135381084
$ws.Range("Q5").Value = 8.839686647053
$ws.Range("R5").Value = 79.557179823477
$ws.Range("S5").Value = 0.0121482549457968
$ws.Range("T5").Value = 0.01214825494579679
# Row 6
$ws.Range("I6").Value = 0.2895605203232839
$ws.Range("J6").Value = 0.2895605203232838
$ws.Range("M6").Value = 18.483507
$ws.Range("N6").Value = 55.450521
$ws.Range("O6").Value = 0.1196447515092806
$ws.Range("P6").Value = 0.1196447515092806
$ws.Range("Q6").Value = 25.209020596893
$ws.Range("R6").Value = 226.881185372037
$ws.Range("S6").Value = 0.0346443965009773
$ws.Range("T6").Value = 0.0346443965009773
# Row 7
$ws.Range("G7").Value = 1.444951
$ws.Range("H7").Value = 4.334853
$ws.Range("I7").Value = 0.3067756404662893
$ws.Range("J7").Value = 0.3067756404662893
$ws.Range("M7").Value = 51.15371566666666
$ws.Range("N7").Value = 153.461147
$ws.Range("O7").Value = 0.3311207986511828
$ws.Range("P7").Value = 0.3311207986511828
$ws.Range("Q7").Value = 73.91461260626565
$ws.Range("R7").Value = 665.2315134563909
$ws.Range("S7").Value = 0.1015797950779258
$ws.Range("T7").Value = 0.1015797950779258
# Row 8
$ws.Range("G8").Value = 1.444951
$ws.Range("H8").Value = 4.334853
$ws.Range("I8").Value = 0.3067756404662893
$ws.Range("J8").Value = 0.3067756404662893
$ws.Range("M8").Value = 53.36146666666667
$ws.Range("O8").Value = 0.3454116915964105
$ws.Range("P8").Value = 0.3454116915964106
$ws.Range("Q8").Value = 77.10470462146667
$ws.Range("R8").Value = 693.9423415932
$ws.Range("S8").Value = 0.1059638929140332
$ws.Range("T8").Value = 0.1059638929140332
# Row 9
$ws.Range("G9").Value = 1.444951
$ws.Range("H9").Value = 4.334853
$ws.Range("I9").Value = 0.3067756404662893
$ws.Range("J9").Value = 0.3067756404662893
$ws.Range("M9").Value = 25.00653133333333
$ws.Range("N9").Value = 75.019594
$ws.Range("O9").Value = 0.1618686447050176
$ws.Range("P9").Value = 0.1618686447050176
$ws.Range("Q9").Value = 36.13321245663133
$ws.Range("R9").Value = 325.198912109682
$ws.Range("S9").Value = 0.04965735715079201
$ws.Range("T9").Value = 0.04965735715079201
# Row 10
$ws.Range("G10").Value = 1.444951
$ws.Range("H10").Value = 4.334853
$ws.Range("I10").Value = 0.3067756404662893
$ws.Range("J10").Value = 0.3067756404662893
$ws.Range("M10").Value = 6.481347
$ws.Range("N10").Value = 19.444041
$ws.Range("O10").Value = 0.0419541135381084
$ws.Range("P10").Value = 0.0419541135381084
$ws.Range("Q10").Value = 9.365228828996999
$ws.Range("R10").Value = 84.287059460973
$ws.Range("S10").Value = 0.01287050005084862
$ws.Range("T10").Value = 0.01287050005084862
# Row 11
$ws.Range("G11").Value = 1.444951
$ws.Range("H11").Value = 4.334853
$ws.Range("I11").Value = 0.3067756404662893
$ws.Range("J11").Value = 0.3067756404662893
$ws.Range("M11").Value = 18.483507
$ws.Range("N11").Value = 55.450521
$ws.Range("O11").Value = 0.1196447515092806
$ws.Range("P11").Value = 0.1196447515092806
$ws.Range("Q11").Value = 26.707761923157
$ws.Range("R11").Value = 240.369857308413
$ws.Range("S11").Value = 0.03670409527268959
$ws.Range("T11").Value = 0.0367040952726896
# Row 12
$ws.Range("G12").Value = 0.680678
$ws.Range("H12").Value = 2.042034
$ws.Range("I12").Value = 0.1445138481521608
$ws.Range("J12").Value = 0.1445138481521608
$ws.Range("M12").Value = 51.15371566666666
$ws.Range("N12").Value = 153.461147
$ws.Range("O12").Value = 0.3311207986511828
$ws.Range("P12").Value = 0.3311207986511828
$ws.Range("Q12").Value = 34.81920887255533
$ws.Range("R12").Value = 313.372879852998
$ws.Range("S12").Value = 0.04785154081629925
$ws.Range("T12").Value = 0.04785154081629924
# Row 13
$ws.Range("G13").Value = 0.680678
$ws.Range("H13").Value = 2.042034
$ws.Range("I13").Value = 0.1445138481521608
$ws.Range("J13").Value = 0.1445138481521608
$ws.Range("M13").Value = 53.36146666666667
$ws.Range("O13").Value = 0.3454116915964105
$ws.Range("P13").Value = 0.3454116915964106
$ws.Range("Q13").Value = 36.32197640773334
$ws.Range("R13").Value = 326.8977876696001
$ws.Range("S13").Value = 0.04991677274934467
$ws.Range("T13").Value = 0.04991677274934467
# Row 14
$ws.Range("G14").Value = 0.680678
$ws.Range("H14").Value = 2.042034
$ws.Range("I14").Value = 0.1445138481521608
$ws.Range("J14").Value = 0.1445138481521608
$ws.Range("M14").Value = 25.00653133333333
$ws.Range("N14").Value = 75.019594
$ws.Range("O14").Value = 0.1618686447050176
$ws.Range("P14").Value = 0.1618686447050176
$ws.Range("Q14").Value = 17.02139573491067
$ws.Range("R14").Value = 153.192561614196
$ws.Range("S14").Value = 0.02339226074149699
$ws.Range("T14").Value = 0.02339226074149698
# Row 15
$ws.Range("G15").Value = 0.680678
$ws.Range("H15").Value = 2.042034
$ws.Range("I15").Value = 0.1445138481521608
$ws.Range("J15").Value = 0.1445138481521608
$ws.Range("M15").Value = 6.481347
$ws.Range("N15").Value = 19.444041
$ws.Range("O15").Value = 0.0419541135381084
$ws.Range("P15").Value = 0.0419541135381084
$ws.Range("Q15").Value = 4.411710313266
$ws.Range("R15").Value = 39.705392819394
$ws.Range("S15").Value = 0.006062950393204711
$ws.Range("T15").Value = 0.00606295039320471
# Row 16
$ws.Range("G16").Value = 0.680678
$ws.Range("H16").Value = 2.042034
$ws.Range("I16").Value = 0.1445138481521608
$ws.Range("J16").Value = 0.1445138481521608
$ws.Range("M16").Value = 18.483507
$ws.Range("N16").Value = 55.450521
$ws.Range("O16").Value = 0.1196447515092806
$ws.Range("P16").Value = 0.1196447515092806
$ws.Range("Q16").Value = 12.581316577746
$ws.Range("R16").Value = 113.231849199714
$ws.Range("S16").Value = 0.01729032345181519
$ws.Range("T16").Value = 0.01729032345181519
# Row 17
$ws.Range("G17").Value = 1.220628333333333
$ws.Range("H17").Value = 3.661885
$ws.Range("I17").Value = 0.2591499910582661
$ws.Range("J17").Value = 0.2591499910582661
$ws.Range("M17").Value = 51.15371566666666
$ws.Range("N17").Value = 153.461147
$ws.Range("O17").Value = 0.3311207986511828
$ws.Range("P17").Value = 0.3311207986511828
$ws.Range("Q17").Value = 62.43967469801056
$ws.Range("R17").Value = 561.9570722820949
$ws.Range("S17").Value = 0.08580995200965998
$ws.Range("T17").Value = 0.08580995200965995
# Row 18
$ws.Range("G18").Value = 1.220628333333333
$ws.Range("H18").Value = 3.661885
$ws.Range("I18").Value = 0.2591499910582661
$ws.Range("J18").Value = 0.2591499910582661
$ws.Range("M18").Value = 53.36146666666667
$ws.Range("O18").Value = 0.3454116915964105
$ws.Range("P18").Value = 0.3454116915964106
$ws.Range("Q18").Value = 65.13451812155557
$ws.Range("R18").Value = 586.2106630940001
$ws.Range("S18").Value = 0.08951343678863036
$ws.Range("T18").Value = 0.08951343678863036
# Row 19
$ws.Range("G19").Value = 1.220628333333333
$ws.Range("H19").Value = 3.661885
$ws.Range("I19").Value = 0.2591499910582661
$ws.Range("J19").Value = 0.2591499910582661
$ws.Range("M19").Value = 25.00653133333333
$ws.Range("N19").Value = 75.019594
$ws.Range("O19").Value = 0.1618686447050176
$ws.Range("P19").Value = 0.1618686447050176
$ws.Range("Q19").Value = 30.52368066385445
$ws.Range("R19").Value = 274.71312597469
$ws.Range("S19").Value = 0.04194825782791898
$ws.Range("T19").Value = 0.04194825782791897
# Row 20
$ws.Range("G20").Value = 1.220628333333333
$ws.Range("H20").Value = 3.661885
$ws.Range("I20").Value = 0.2591499910582661
$ws.Range("J20").Value = 0.2591499910582661
$ws.Range("M20").Value = 6.481347
$ws.Range("N20").Value = 19.444041
$ws.Range("O20").Value = 0.0419541135381084
$ws.Range("P20").Value = 0.0419541135381084
$ws.Range("Q20").Value = 7.911315786365
$ws.Range("R20").Value = 71.20184207728499
$ws.Range("S20").Value = 0.01087240814825827
$ws.Range("T20").Value = 0.01087240814825827
# Row 21
$ws.Range("G21").Value = 1.220628333333333
$ws.Range("H21").Value = 3.661885
$ws.Range("I21").Value = 0.2591499910582661
$ws.Range("J21").Value = 0.2591499910582661
$ws.Range("M21").Value = 18.483507
$ws.Range("N21").Value = 55.450521
$ws.Range("O21").Value = 0.1196447515092806
$ws.Range("P21").Value = 0.1196447515092806
$ws.Range("Q21").Value = 22.561492343565
$ws.Range("R21").Value = 240.369857308413
$ws.Range("S21").Value = 0.03670409527268959
$ws.Range("T21").Value = 0.0367040952726896
